$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $val) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value2 = $val
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '257.30'
Set-TextValue 'E2' '-0.50%'
Set-TextValue 'E3' '0.34%'
Set-TextValue 'D4' '4.582'
Set-TextValue 'E4' '-6.19%'
Set-TextValue 'D5' '0.05895'
Set-TextValue 'E5' '-1.07%'
Set-TextValue 'D6' '6.629'
Set-TextValue 'D7' '0.8514'
Set-TextValue 'E7' '-2.69%'
Set-TextValue 'D8' '0.9429'
Set-TextValue 'E8' '-2.15%'
Set-TextValue 'B9' 'One'
Set-TextValue 'C9' 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
Set-TextValue 'D9' '0.01046'
Set-TextValue 'E9' '1,621.74%'
Set-TextValue 'B10' 'WazirX'
Set-TextValue 'C10' 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
Set-TextValue 'D10' '0.1389'
Set-TextValue 'E10' '-1.86%'
Set-TextValue 'B11' 'LiechtensteinCryptoassetsExchange'
Set-TextValue 'C11' 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
Set-TextValue 'D11' '0.04892'
Set-TextValue 'E11' '36.19%'
Set-TextValue 'B12' 'MandalaExchangeToken'
Set-TextValue 'C12' 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
Set-TextValue 'D12' '0.07074'
Set-TextValue 'E12' '-1.45%'
Set-TextValue 'B13' 'BitrueCoin'
Set-TextValue 'C13' 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
Set-TextValue 'D13' '0.03072'
Set-TextValue 'E13' '-2.04%'
Set-TextValue 'B14' 'BitMartToken'
Set-TextValue 'C14' 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
Set-TextValue 'D14' '0.09134'
Set-TextValue 'E14' '-1.07%'
Set-TextValue 'B15' 'BitForexToken'
Set-TextValue 'C15' 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
Set-TextValue 'D15' '0.001532'
Set-TextValue 'E15' '-1.17%'
Set-TextValue 'D16' '0.006030'
Set-TextValue 'E16' '0.55%'
Set-TextValue 'D17' '3.496'
Set-TextValue 'E17' '0.35%'
Set-TextValue 'D18' '3.182'
Set-TextValue 'E20' '-2.90%'
Set-TextValue 'E21' '-2.78%'
Set-TextValue 'D22' '3.961'
Set-TextValue 'E22' '12.41%'
Set-TextValue 'D23' '0.04253'
Set-TextValue 'E23' '1.00%'
Set-TextValue 'D24' '0.001222'
Set-TextValue 'E24' '0.30%'
Set-TextValue 'E25' '-5.19%'
Set-TextValue 'E26' '0.09%'
Set-TextValue 'D27' '0.0001524'
Set-TextValue 'E27' '2.13%'
Set-TextValue 'E40' '-0.38%'
Set-TextValue 'D41' '0.006232'
Set-TextValue 'E41' '5.87%'
Set-TextValue 'D42' '0.1101'
Set-TextValue 'E42' '-0.20%'
Set-TextValue 'D43' '0.002201'
Set-TextValue 'E43' '0.10%'
Set-TextValue 'D44' '0.01418'
Set-TextValue 'E44' '35.17%'
Set-TextValue 'D45' '0.00005380'
Set-TextValue 'E45' '-2.02%'
Set-TextValue 'E46' '0.10%'
Set-TextValue 'D47' '0.06690'
Set-TextValue 'E47' '-38.64%'
Set-TextValue 'E48' '11,584.52%'
Set-TextValue 'D49' '0.00002101'
Set-TextValue 'E49' '0.10%'
Set-TextValue 'D50' '0.0002001'
Set-TextValue 'E50' '0.10%'
